$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6 (item id 4564)
$ws.Cells.Item(6, 8).Value = 83  # H6: 88.75 -> 83
$ws.Cells.Item(6, 9).Value = 83  # I6: 88.75 -> 83
$ws.Cells.Item(6, 11).Value = 249  # K6: 266.25 -> 249
$ws.Cells.Item(6, 13).Value = -137  # M6: -154.25 -> -137

# Row 18 (item id 5471)
$ws.Cells.Item(18, 8).Value = 1254.3334  # H18: 760.7143 -> 1254.3334
$ws.Cells.Item(18, 9).Value = 567.5714  # I18: 720.8333 -> 567.5714
$ws.Cells.Item(18, 10).Value = 3658  # J18: 1000 -> 3658
$ws.Cells.Item(18, 11).Value = 567.5714  # K18: 720.8333 -> 567.5714
$ws.Cells.Item(18, 12).Value = 3658  # L18: 1000 -> 3658
$ws.Cells.Item(18, 13).Value = -283.5714  # M18: -436.8333 -> -283.5714
$ws.Cells.Item(18, 14).Value = -4226  # N18: -1568 -> -4226

# Row 33 (item id 5512)
$ws.Cells.Item(33, 8).Value = 730.6  # H33: 693.375 -> 730.6
$ws.Cells.Item(33, 10).Value = 1605.6  # J33: 1360.5 -> 1605.6
$ws.Cells.Item(33, 12).Value = 1605.6  # L33: 1360.5 -> 1605.6
$ws.Cells.Item(33, 14).Value = -2063.6  # N33: -1818.5 -> -2063.6

# Row 40 (item id 5505)
$ws.Cells.Item(40, 8).Value = 2450  # H40: 1967 -> 2450
$ws.Cells.Item(40, 10).Value = 2450  # J40: 1967 -> 2450
$ws.Cells.Item(40, 12).Value = 2450  # L40: 1967 -> 2450
$ws.Cells.Item(40, 14).Value = -2800  # N40: -2317 -> -2800

# Row 46 (item id 4584)
$ws.Cells.Item(46, 8).Value = 8900  # H46: 8993.5 -> 8900
$ws.Cells.Item(46, 10).Value = 5800  # J46: 5987 -> 5800
$ws.Cells.Item(46, 12).Value = 17400  # L46: 17961 -> 17400
$ws.Cells.Item(46, 14).Value = -17638  # N46: -18199 -> -17638

# Row 60 (item id 4584)
$ws.Cells.Item(60, 8).Value = 8900  # H60: 8993.5 -> 8900
$ws.Cells.Item(60, 10).Value = 5800  # J60: 5987 -> 5800
$ws.Cells.Item(60, 12).Value = 17400  # L60: 17961 -> 17400
$ws.Cells.Item(60, 14).Value = -18368  # N60: -18929 -> -18368

# Row 69 (item id 12616)
$ws.Cells.Item(69, 8).Value = 10000  # H69: 5250 -> 10000
$ws.Cells.Item(69, 10).Value = 10000  # J69: 5250 -> 10000
$ws.Cells.Item(69, 12).Value = 30000  # L69: 15750 -> 30000
$ws.Cells.Item(69, 14).Value = -31748  # N69: -17498 -> -31748

# Row 70 (item id 12604)
$ws.Cells.Item(70, 8).Value = 1812.5  # H70: 1833.3334 -> 1812.5
$ws.Cells.Item(70, 10).Value = 1812.5  # J70: 1833.3334 -> 1812.5
$ws.Cells.Item(70, 12).Value = 5437.5  # L70: 5500.0002 -> 5437.5
$ws.Cells.Item(70, 14).Value = -5977.5  # N70: -6040.0002 -> -5977.5

# Row 72 (item id 12616)
$ws.Cells.Item(72, 8).Value = 10000  # H72: 5250 -> 10000
$ws.Cells.Item(72, 10).Value = 10000  # J72: 5250 -> 10000
$ws.Cells.Item(72, 12).Value = 90000  # L72: 47250 -> 90000
$ws.Cells.Item(72, 14).Value = -98736  # N72: -55986 -> -98736

# Row 73 (item id 12604)
$ws.Cells.Item(73, 8).Value = 1812.5  # H73: 1833.3334 -> 1812.5
$ws.Cells.Item(73, 10).Value = 1812.5  # J73: 1833.3334 -> 1812.5
$ws.Cells.Item(73, 12).Value = 5437.5  # L73: 5500.0002 -> 5437.5
$ws.Cells.Item(73, 14).Value = -7309.5  # N73: -7372.0002 -> -7309.5

# Row 113 (item id 27775)
$ws.Cells.Item(113, 8).Value = 2000  # H113: 1857.5385 -> 2000
$ws.Cells.Item(113, 9).Value = 2000  # I113: 1845.6666 -> 2000
$ws.Cells.Item(113, 11).Value = 2000  # K113: 1845.6666 -> 2000
$ws.Cells.Item(113, 13).Value = 1254  # M113: 1408.3334 -> 1254

# Row 137 (item id 44013)
$ws.Cells.Item(137, 8).Value = 4999  # H137: 1949.5 -> 4999
$ws.Cells.Item(137, 9).Value = 0  # I137: 1949.5 -> 0
$ws.Cells.Item(137, 10).Value = 4999  # J137: 0 -> 4999
$ws.Cells.Item(137, 11).Value = 0  # K137: 5848.5 -> 0
$ws.Cells.Item(137, 12).Value = 14997  # L137: 0 -> 14997
$ws.Cells.Item(137, 13).ClearContents()  # M137: -3298.5 -> (removed)
$ws.Cells.Item(137, 14).Value = -20097  # N137: None -> -20097

$ws = $wb.Worksheets.Item("ARM")
# Row 63 (item id 12528)
$ws.Cells.Item(63, 8).Value = 0  # H63: 320 -> 0
$ws.Cells.Item(63, 9).Value = 0  # I63: 320 -> 0
$ws.Cells.Item(63, 11).Value = 0  # K63: 320 -> 0
$ws.Cells.Item(63, 13).ClearContents()  # M63: 366 -> (removed)

# Row 66 (item id 12528)
$ws.Cells.Item(66, 8).Value = 0  # H66: 320 -> 0
$ws.Cells.Item(66, 9).Value = 0  # I66: 320 -> 0
$ws.Cells.Item(66, 11).Value = 0  # K66: 1600 -> 0
$ws.Cells.Item(66, 13).ClearContents()  # M66: 1832 -> (removed)

# Row 97 (item id 19941)
$ws.Cells.Item(97, 8).Value = 929.5  # H97: 0 -> 929.5
$ws.Cells.Item(97, 9).Value = 904.6667  # I97: 0 -> 904.6667
$ws.Cells.Item(97, 10).Value = 1004  # J97: 0 -> 1004
$ws.Cells.Item(97, 11).Value = 904.6667  # K97: 0 -> 904.6667
$ws.Cells.Item(97, 12).Value = 1004  # L97: 0 -> 1004
$ws.Cells.Item(97, 13).Value = -408.6667  # M97: None -> -408.6667
$ws.Cells.Item(97, 14).Value = -1996  # N97: None -> -1996

# Row 132 (item id 43997)
$ws.Cells.Item(132, 8).Value = 3006  # H132: 5012 -> 3006
$ws.Cells.Item(132, 9).Value = 3006  # I132: 5012 -> 3006
$ws.Cells.Item(132, 11).Value = 9018  # K132: 15036 -> 9018
$ws.Cells.Item(132, 13).Value = -6488  # M132: -12506 -> -6488

$ws = $wb.Worksheets.Item("BSM")
# Row 5 (item id 1750)
$ws.Cells.Item(5, 8).Value = 617.8333  # H5: 694.6 -> 617.8333
$ws.Cells.Item(5, 9).Value = 379.33334  # I5: 401.33334 -> 379.33334
$ws.Cells.Item(5, 10).Value = 856.3333  # J5: 1134.5 -> 856.3333
$ws.Cells.Item(5, 11).Value = 379.33334  # K5: 401.33334 -> 379.33334
$ws.Cells.Item(5, 12).Value = 856.3333  # L5: 1134.5 -> 856.3333
$ws.Cells.Item(5, 13).Value = -266.33334  # M5: -288.33334 -> -266.33334
$ws.Cells.Item(5, 14).Value = -1082.3333  # N5: -1360.5 -> -1082.3333

# Row 11 (item id 2481)
$ws.Cells.Item(11, 8).Value = 3945.8  # H11: 5883 -> 3945.8
$ws.Cells.Item(11, 9).Value = 6052  # I11: 12004 -> 6052
$ws.Cells.Item(11, 10).Value = 2541.6667  # J11: 2822.5 -> 2541.6667
$ws.Cells.Item(11, 11).Value = 6052  # K11: 12004 -> 6052
$ws.Cells.Item(11, 12).Value = 2541.6667  # L11: 2822.5 -> 2541.6667
$ws.Cells.Item(11, 13).Value = -5912  # M11: -11864 -> -5912
$ws.Cells.Item(11, 14).Value = -2821.6667  # N11: -3102.5 -> -2821.6667

# Row 94 (item id 19939)
$ws.Cells.Item(94, 8).Value = 625  # H94: 400 -> 625
$ws.Cells.Item(94, 9).Value = 566.6667  # I94: 400 -> 566.6667
$ws.Cells.Item(94, 10).Value = 800  # J94: 0 -> 800
$ws.Cells.Item(94, 11).Value = 566.6667  # K94: 400 -> 566.6667
$ws.Cells.Item(94, 12).Value = 800  # L94: 0 -> 800
$ws.Cells.Item(94, 13).Value = -115.6667  # M94: 51 -> -115.6667
$ws.Cells.Item(94, 14).Value = -1702  # N94: None -> -1702

# Row 105 (item id 19947)
$ws.Cells.Item(105, 8).Value = 73061.8  # H105: 61068.168 -> 73061.8
$ws.Cells.Item(105, 9).Value = 53837  # I105: 40652.75 -> 53837
$ws.Cells.Item(105, 11).Value = 53837  # K105: 40652.75 -> 53837
$ws.Cells.Item(105, 13).Value = -52090  # M105: -38905.75 -> -52090

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (item id 5361)
$ws.Cells.Item(7, 8).Value = 191.66667  # H7: 187.5 -> 191.66667
$ws.Cells.Item(7, 9).Value = 75  # I7: 83.333336 -> 75
$ws.Cells.Item(7, 10).Value = 425  # J7: 500 -> 425
$ws.Cells.Item(7, 11).Value = 75  # K7: 83.333336 -> 75
$ws.Cells.Item(7, 12).Value = 425  # L7: 500 -> 425
$ws.Cells.Item(7, 13).Value = 38  # M7: 29.666664 -> 38
$ws.Cells.Item(7, 14).Value = -651  # N7: -726 -> -651

# Row 31 (item id 44023)
$ws.Cells.Item(31, 8).Value = 3797.4  # H31: 3799.4 -> 3797.4
$ws.Cells.Item(31, 9).Value = 1995  # I31: 2000 -> 1995
$ws.Cells.Item(31, 11).Value = 1995  # K31: 2000 -> 1995
$ws.Cells.Item(31, 13).Value = -1700  # M31: -1705 -> -1700

# Row 34 (item id 44023)
$ws.Cells.Item(34, 8).Value = 3797.4  # H34: 3799.4 -> 3797.4
$ws.Cells.Item(34, 9).Value = 1995  # I34: 2000 -> 1995
$ws.Cells.Item(34, 11).Value = 1995  # K34: 2000 -> 1995
$ws.Cells.Item(34, 13).Value = -1793  # M34: -1798 -> -1793

# Row 47 (item id 1920)
$ws.Cells.Item(47, 8).Value = 0  # H47: 75000 -> 0
$ws.Cells.Item(47, 10).Value = 0  # J47: 75000 -> 0
$ws.Cells.Item(47, 12).Value = 0  # L47: 75000 -> 0
$ws.Cells.Item(47, 14).ClearContents()  # N47: -76132 -> (removed)

# Row 50 (item id 1862)
$ws.Cells.Item(50, 8).Value = 0  # H50: 20000 -> 0
$ws.Cells.Item(50, 10).Value = 0  # J50: 20000 -> 0
$ws.Cells.Item(50, 12).Value = 0  # L50: 20000 -> 0
$ws.Cells.Item(50, 14).ClearContents()  # N50: -21250 -> (removed)

# Row 60 (item id 1937)
$ws.Cells.Item(60, 8).Value = 0  # H60: 2000 -> 0
$ws.Cells.Item(60, 9).Value = 0  # I60: 2000 -> 0
$ws.Cells.Item(60, 11).Value = 0  # K60: 2000 -> 0
$ws.Cells.Item(60, 13).ClearContents()  # M60: -1489 -> (removed)

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (item id 4847)
$ws.Cells.Item(2, 8).Value = 2.4  # H2: 2.1666667 -> 2.4
$ws.Cells.Item(2, 9).Value = 3  # I2: 2.3333333 -> 3
$ws.Cells.Item(2, 11).Value = 18  # K2: 13.9999998 -> 18
$ws.Cells.Item(2, 13).Value = 95  # M2: 99.0000002 -> 95

# Row 7 (item id 4728)
$ws.Cells.Item(7, 8).Value = 50  # H7: 0 -> 50
$ws.Cells.Item(7, 9).Value = 50  # I7: 0 -> 50
$ws.Cells.Item(7, 11).Value = 150  # K7: 0 -> 150
$ws.Cells.Item(7, 13).Value = -38  # M7: None -> -38

# Row 8 (item id 16734)
$ws.Cells.Item(8, 8).Value = 85  # H8: 81.666664 -> 85
$ws.Cells.Item(8, 9).Value = 85  # I8: 81.666664 -> 85
$ws.Cells.Item(8, 11).Value = 255  # K8: 244.999992 -> 255
$ws.Cells.Item(8, 13).Value = -116  # M8: -105.999992 -> -116

# Row 12 (item id 4854)
$ws.Cells.Item(12, 8).Value = 177.33333  # H12: 190.66667 -> 177.33333
$ws.Cells.Item(12, 10).Value = 177.33333  # J12: 190.66667 -> 177.33333
$ws.Cells.Item(12, 12).Value = 531.99999  # L12: 572.00001 -> 531.99999
$ws.Cells.Item(12, 14).Value = -877.99999  # N12: -918.00001 -> -877.99999

# Row 13 (item id 4657)
$ws.Cells.Item(13, 8).Value = 293.25  # H13: 238 -> 293.25
$ws.Cells.Item(13, 9).Value = 307  # I13: 238 -> 307
$ws.Cells.Item(13, 10).Value = 252  # J13: 0 -> 252
$ws.Cells.Item(13, 11).Value = 921  # K13: 714 -> 921
$ws.Cells.Item(13, 12).Value = 756  # L13: 0 -> 756
$ws.Cells.Item(13, 13).Value = -753  # M13: -546 -> -753
$ws.Cells.Item(13, 14).Value = -1092  # N13: None -> -1092

# Row 34 (item id 4749)
$ws.Cells.Item(34, 8).Value = 6155.8  # H34: 4529.8887 -> 6155.8
$ws.Cells.Item(34, 9).Value = 639.5  # I34: 592.25 -> 639.5
$ws.Cells.Item(34, 10).Value = 9833.333000000001  # J34: 7680 -> 9833.333000000001
$ws.Cells.Item(34, 11).Value = 1918.5  # K34: 1776.75 -> 1918.5
$ws.Cells.Item(34, 12).Value = 29499.999  # L34: 23040 -> 29499.999
$ws.Cells.Item(34, 13).Value = -1834.5  # M34: -1692.75 -> -1834.5
$ws.Cells.Item(34, 14).Value = -29667.999  # N34: -23208 -> -29667.999

# Row 39 (item id 4712)
$ws.Cells.Item(39, 8).Value = 10083  # H39: 7921.4287 -> 10083
$ws.Cells.Item(39, 10).Value = 10083  # J39: 7921.4287 -> 10083
$ws.Cells.Item(39, 12).Value = 30249  # L39: 23764.2861 -> 30249
$ws.Cells.Item(39, 14).Value = -30837  # N39: -24352.2861 -> -30837

# Row 55 (item id 4733)
$ws.Cells.Item(55, 8).Value = 6699.8887  # H55: 7489.9 -> 6699.8887
$ws.Cells.Item(55, 9).Value = 1650  # I55: 1700 -> 1650
$ws.Cells.Item(55, 10).Value = 8142.7144  # J55: 8937.375 -> 8142.7144
$ws.Cells.Item(55, 11).Value = 4950  # K55: 5100 -> 4950
$ws.Cells.Item(55, 12).Value = 24428.1432  # L55: 26812.125 -> 24428.1432
$ws.Cells.Item(55, 13).Value = -4773  # M55: -4923 -> -4773
$ws.Cells.Item(55, 14).Value = -24782.1432  # N55: -27166.125 -> -24782.1432

# Row 88 (item id 12851)
$ws.Cells.Item(88, 8).Value = 0  # H88: 15000 -> 0
$ws.Cells.Item(88, 10).Value = 0  # J88: 15000 -> 0
$ws.Cells.Item(88, 12).Value = 0  # L88: 45000 -> 0
$ws.Cells.Item(88, 14).ClearContents()  # N88: -45856 -> (removed)

# Row 91 (item id 12851)
$ws.Cells.Item(91, 8).Value = 0  # H91: 15000 -> 0
$ws.Cells.Item(91, 10).Value = 0  # J91: 15000 -> 0
$ws.Cells.Item(91, 12).Value = 0  # L91: 45000 -> 0
$ws.Cells.Item(91, 14).ClearContents()  # N91: -47964 -> (removed)

$ws = $wb.Worksheets.Item("GSM")
# Row 42 (item id 27213)
$ws.Cells.Item(42, 8).Value = 0  # H42: 60000 -> 0
$ws.Cells.Item(42, 10).Value = 0  # J42: 60000 -> 0
$ws.Cells.Item(42, 12).Value = 0  # L42: 60000 -> 0
$ws.Cells.Item(42, 14).ClearContents()  # N42: -60970 -> (removed)

# Row 97 (item id 19940)
$ws.Cells.Item(97, 8).Value = 200  # H97: 156.33333 -> 200
$ws.Cells.Item(97, 9).Value = 200  # I97: 134.5 -> 200
$ws.Cells.Item(97, 10).Value = 0  # J97: 200 -> 0
$ws.Cells.Item(97, 11).Value = 200  # K97: 134.5 -> 200
$ws.Cells.Item(97, 12).Value = 0  # L97: 200 -> 0
$ws.Cells.Item(97, 13).Value = 296  # M97: 361.5 -> 296
$ws.Cells.Item(97, 14).ClearContents()  # N97: -1192 -> (removed)

# Row 115 (item id 27213)
$ws.Cells.Item(115, 8).Value = 0  # H115: 60000 -> 0
$ws.Cells.Item(115, 10).Value = 0  # J115: 60000 -> 0
$ws.Cells.Item(115, 12).Value = 0  # L115: 60000 -> 0
$ws.Cells.Item(115, 14).ClearContents()  # N115: -62350 -> (removed)

$ws = $wb.Worksheets.Item("LTW")
# Row 94 (item id 18067)
$ws.Cells.Item(94, 8).Value = 980000  # H94: 95000 -> 980000
$ws.Cells.Item(94, 10).Value = 980000  # J94: 95000 -> 980000
$ws.Cells.Item(94, 12).Value = 980000  # L94: 95000 -> 980000
$ws.Cells.Item(94, 14).Value = -981352  # N94: -96352 -> -981352

# Row 136 (item id 44060)
$ws.Cells.Item(136, 8).Value = 27875  # H136: 36666.668 -> 27875
$ws.Cells.Item(136, 10).Value = 23250  # J136: 45000 -> 23250
$ws.Cells.Item(136, 12).Value = 69750  # L136: 135000 -> 69750
$ws.Cells.Item(136, 14).Value = -74850  # N136: -140100 -> -74850

$ws = $wb.Worksheets.Item("WVR")
# Row 18 (item id 3543)
$ws.Cells.Item(18, 8).Value = 0  # H18: 100 -> 0
$ws.Cells.Item(18, 9).Value = 0  # I18: 100 -> 0
$ws.Cells.Item(18, 11).Value = 0  # K18: 100 -> 0
$ws.Cells.Item(18, 13).ClearContents()  # M18: 73 -> (removed)

# Row 81 (item id 12596)
$ws.Cells.Item(81, 8).Value = 5597.6  # H81: 6747 -> 5597.6
$ws.Cells.Item(81, 10).Value = 1000  # J81: 0 -> 1000
$ws.Cells.Item(81, 12).Value = 2000  # L81: 0 -> 2000
$ws.Cells.Item(81, 14).Value = -4122  # N81: None -> -4122

# Row 84 (item id 12596)
$ws.Cells.Item(84, 8).Value = 5597.6  # H84: 6747 -> 5597.6
$ws.Cells.Item(84, 10).Value = 1000  # J84: 0 -> 1000
$ws.Cells.Item(84, 12).Value = 10000  # L84: 0 -> 10000
$ws.Cells.Item(84, 14).Value = -20608  # N84: None -> -20608

# Row 122 (item id 36208)
$ws.Cells.Item(122, 8).Value = 1708.8  # H122: 1924 -> 1708.8
$ws.Cells.Item(122, 9).Value = 1708.8  # I122: 1924 -> 1708.8
$ws.Cells.Item(122, 11).Value = 5126.4  # K122: 5772 -> 5126.4
$ws.Cells.Item(122, 13).Value = -2676.4  # M122: -3322 -> -2676.4
